# Remove the "target_ids" column from the target_data worksheet / table.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("target_data")

$table = $ws.ListObjects.Item("Table4")
$col = $table.ListColumns.Item("target_ids")
$col.Delete()
